$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = 281
$ws.Range("C17").Value = 121
$ws.Range("D17").Value = 24

$ws.Range("B18").Value = 1334
$ws.Range("C18").Value = 677
$ws.Range("D18").Value = 113

$ws.Range("B19").Value = 1343
$ws.Range("C19").Value = 865
$ws.Range("D19").Value = 274

$ws.Range("B20").Value = 303
$ws.Range("C20").Value = 249
$ws.Range("D20").Value = 29

$ws.Range("B21").Value = 367
$ws.Range("C21").Value = 129
$ws.Range("D21").Value = 62

$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 444
$ws.Range("D22").Value = 109
